$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# The shape-id allocator in this host assigns the next *unused* id (starting
# at 2) to every newly created shape, independent of the ids already present
# in the slide's XML. The target OOXML expects the new wrapper group to be
# born with id="16" / name="Group 15" (matching real PowerPoint's "Group N"
# -> id N+1 convention). Burn through the intervening ids with disposable
# textboxes (created, then removed) so the group we keep lands on id 16.
$burn = @()
for ($i = 0; $i -lt 9; $i++) {
    $burn += $s.Shapes.AddTextbox(1, 0, 0, 1, 1)
}

# Group the five existing top-level shapes (picture, dashed line, the two
# "max value" / arrow shapes, and the "threshold value" textbox) together,
# preserving their current z-order.
$range = $s.Shapes.Range(@(1, 2, 3, 4, 5))
$group = $range.Group()
$group.Name = "Group 15"

foreach ($tmp in $burn) {
    $tmp.Delete()
}
